# Locate the target paragraph: "диаметр малого ступенчатого отверстия крышки"
$d = $word.ActiveDocument

$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*диаметр малого ступенчатого отверстия крышки*") {
        $target = $p
        break
    }
}

# Update indentation: left 567 -> 0 twips, firstLine 709 -> 1276 twips
# (Word's ParagraphFormat works in points; 20 twips = 1 point)
$target.Format.LeftIndent = 0
$target.Format.FirstLineIndent = 1276 / 20

# Move the "_GoBack" bookmark to the very start of this paragraph (before its
# first run). Word keeps only one bookmark per name, so adding it here
# removes the previous one automatically.
$startRange = $d.Range($target.Range.Start, $target.Range.Start)
$d.Bookmarks.Add("_GoBack", $startRange)
